# Update the LR-pairs sheet (Gnai2-Adora1) with the new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster: ECs, Target cluster -> ECs)
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02215633333333333
$ws.Range("N2").Value = 0.066469
$ws.Range("O2").Value = 0.03596233285271019
$ws.Range("P2").Value = 0.03596233285271019
$ws.Range("Q2").Value = 3.764707935044334
$ws.Range("R2").Value = 33.882371415399
$ws.Range("S2").Value = 0.01597153707754774
$ws.Range("T2").Value = 0.01597153707754774

# Row 3 (Sending cluster: ECs, Target cluster -> MuSCs)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5939420000000001
$ws.Range("N3").Value = 1.781826
$ws.Range("O3").Value = 0.9640376671472898
$ws.Range("P3").Value = 0.9640376671472898
$ws.Range("Q3").Value = 100.920045149894
$ws.Range("R3").Value = 908.2804063490461
$ws.Range("S3").Value = 0.4281469560959031
$ws.Range("T3").Value = 0.4281469560959031

# Row 4 (Sending cluster: FAPs, Target cluster -> ECs)
$ws.Range("D4").Value = "ECs"
$ws.Range("I4").Value = 0.1787346690539575
$ws.Range("J4").Value = 0.1787346690539575
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02215633333333333
$ws.Range("N4").Value = 0.066469
$ws.Range("O4").Value = 0.03596233285271019
$ws.Range("P4").Value = 0.03596233285271019
$ws.Range("Q4").Value = 1.515099769989
$ws.Range("R4").Value = 13.635897929901
$ws.Range("S4").Value = 0.00642771566083742
$ws.Range("T4").Value = 0.00642771566083742

# Row 5 (Sending cluster: FAPs, Target cluster -> MuSCs)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5939420000000001
$ws.Range("N5").Value = 1.781826
$ws.Range("O5").Value = 0.9640376671472898
$ws.Range("P5").Value = 0.9640376671472898
$ws.Range("Q5").Value = 40.61508617190601
$ws.Range("R5").Value = 365.535775547154
$ws.Range("S5").Value = 0.1723069533931201
$ws.Range("T5").Value = 0.1723069533931201

# Row 6 (Sending cluster: MuSCs, Target cluster -> ECs)
$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 53.27463399999999
$ws.Range("H6").Value = 159.823902
$ws.Range("I6").Value = 0.1392470275793777
$ws.Range("J6").Value = 0.1392470275793778
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02215633333333333
$ws.Range("N6").Value = 0.066469
$ws.Range("O6").Value = 0.03596233285271019
$ws.Range("P6").Value = 0.03596233285271019
$ws.Range("Q6").Value = 1.180370549115333
$ws.Range("R6").Value = 10.623334942038
$ws.Range("S6").Value = 0.005007647954560098
$ws.Range("T6").Value = 0.005007647954560099

# Row 7 (Sending cluster: MuSCs, Target cluster -> MuSCs)
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 53.27463399999999
$ws.Range("H7").Value = 159.823902
$ws.Range("I7").Value = 0.1392470275793777
$ws.Range("J7").Value = 0.1392470275793778
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5939420000000001
$ws.Range("N7").Value = 1.781826
$ws.Range("O7").Value = 0.9640376671472898
$ws.Range("P7").Value = 0.9640376671472898
$ws.Range("Q7").Value = 31.642042667228
$ws.Range("R7").Value = 284.778384005052
$ws.Range("S7").Value = 0.1342393796248176
$ws.Range("T7").Value = 0.1342393796248177

# Row 8 (Sending cluster: Resolving-Mac, Target cluster -> ECs)
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 91.01828266666666
$ws.Range("H8").Value = 273.054848
$ws.Range("I8").Value = 0.2378998101932138
$ws.Range("J8").Value = 0.2378998101932138
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02215633333333333
$ws.Range("N8").Value = 0.066469
$ws.Range("O8").Value = 0.03596233285271019
$ws.Range("P8").Value = 0.03596233285271019
$ws.Range("Q8").Value = 2.016631410190222
$ws.Range("R8").Value = 18.149682691712
$ws.Range("S8").Value = 0.008555432159764932
$ws.Range("T8").Value = 0.008555432159764932

# Row 9 (Sending cluster: Resolving-Mac, Target cluster -> MuSCs)
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 91.01828266666666
$ws.Range("H9").Value = 273.054848
$ws.Range("I9").Value = 0.2378998101932138
$ws.Range("J9").Value = 0.2378998101932138
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5939420000000001
$ws.Range("N9").Value = 1.781826
$ws.Range("O9").Value = 0.9640376671472898
$ws.Range("P9").Value = 0.9640376671472898
$ws.Range("Q9").Value = 54.05958084360534
$ws.Range("R9").Value = 486.536227592448
$ws.Range("S9").Value = 0.2293443780334489
$ws.Range("T9").Value = 0.2293443780334489
